# #1440 change surveySeries to studySeries
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("study")

# Rename the header labels from surveySeries.de/en to studySeries.de/en
$ws.Range("G1").Value = "studySeries.de"
$ws.Range("H1").Value = "studySeries.en"

# Move the active selection to G2 (as recorded in the saved view state)
$ws.Activate()
$ws.Range("G2").Select() | Out-Null
